$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (Transmitance) values from 1 to 100 for rows 3 through 18
$ws.Range("B3:B18").Value = 100

# Update the selected/active cell to B18
$ws.Range("B18").Select()
